$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$langs = @("ENG","GER","FRE","SPA","ITA","RUS","CHI","UND","POR","POL","CZE","SLO","NOR","JPN","SWE","HRV","DUT","TUR","HUN","UKR","ROM","PERSIAN","GREC","FIN","DAN","ARABE")

for ($i = 0; $i -lt $langs.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = "INSERT INTO LA VALUES ('"
    $ws.Cells.Item($r, 2).Value = $langs[$i] + " "
}

$ws.Rows.Item(1).ClearContents()

$ws.Columns.Item(1).ColumnWidth = 22.67

$ws.Range("A2:C27").Select()
